$d = $word.ActiveDocument

# --- 1. Merge the two runs that were split by the old _GoBack bookmark ---
# "...compare to the DT. Bec" + [bookmark] + "ause the stability..." ->
# "...compare to the DT. Because the stability..." as a single run/bookmark-free span.
$rng = $d.Content
$rng.Find.ClearFormatting()
$oldText = "Bec" + "ause the stability"
$newText = "Because the stability"
$rng.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null

# --- 2. Move the _GoBack bookmark to the empty paragraph right before "Figure #:" ---
# (the paragraph right after "...Due to the minor difference ... in the dataset.")
$anchor = $d.Content
$anchor.Find.ClearFormatting()
$anchor.Find.Execute("Due to the minor difference between the DT, RF and KNN", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

$target = $null
$paras = $d.Paragraphs
for ($i = 1; $i -le $paras.Count; $i++) {
    $p = $paras.Item($i)
    if ($p.Range.Start -le $anchor.Start -and $p.Range.End -ge $anchor.End) {
        $target = $paras.Item($i + 1)
        break
    }
}

$startPos = $target.Range.Start

# Zero-length ranges at an empty paragraph don't resolve reliably, so
# temporarily insert a placeholder character, wrap the bookmark around it,
# then delete the placeholder -- leaving a zero-length bookmark in place.
$insertRange = $d.Range($startPos, $startPos)
$insertRange.InsertBefore("X")

$wrapRange = $d.Range($startPos, $startPos + 1)
$d.Bookmarks.Add("_GoBack", $wrapRange) | Out-Null

$deleteRange = $d.Range($startPos, $startPos + 1)
$deleteRange.Text = ""
